# Insert a new data row at row 46 (pushes existing rows 46-94 down to 47-95)
# and populate it with a new Murcott / Segunda record dated 2022-04-20 (serial 44671).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(46).Insert()

$ws.Cells.Item(46, 1).Value2  = 1
$ws.Cells.Item(46, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(46, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(46, 4).Value2  = 44671
$ws.Cells.Item(46, 5).Value2  = 15
$ws.Cells.Item(46, 6).Value2  = "Fruta"
$ws.Cells.Item(46, 7).Value2  = 100102
$ws.Cells.Item(46, 8).Value2  = "Cítricos"
$ws.Cells.Item(46, 9).Value2  = 100102004
$ws.Cells.Item(46, 10).Value2 = "Mandarina"
$ws.Cells.Item(46, 11).Value2 = "Murcott"
$ws.Cells.Item(46, 12).Value2 = "Segunda"
$ws.Cells.Item(46, 13).Value2 = 300
$ws.Cells.Item(46, 14).Value2 = 19000
$ws.Cells.Item(46, 15).Value2 = 20000
$ws.Cells.Item(46, 16).Value2 = 19500
$ws.Cells.Item(46, 17).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(46, 18).Value2 = "Región de Coquimbo"
$ws.Cells.Item(46, 19).Value2 = 975
$ws.Cells.Item(46, 20).Value2 = 20
